# Updates the voltage-magnitude results table (vm_pu) for the "380 kV" case.
# Columns B-F and I-N for rows 2-25 get new simulation values; columns G/H
# (slack bus / missing bus) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values, one row per data row (rows 2..25), columns in order:
# B, C, D, E, F, I, J, K, L, M, N
$data = @(
    @(1.02,1.017167970187033,1.050404131575146,1.018584458628503,1.053171572303187,1.039191247027457,1.02238362288184,1.053158503322448,1.02143023957654,1.055918277059859,1.023835524895373),
    @(1.02,1.018257849859831,1.051022076573956,1.019513524279183,1.05398574379513,1.039306519346215,1.023108534113564,1.05358938936466,1.02216480748012,1.056545437182579,1.024561465584194),
    @(1.02,1.018963575465803,1.051421598318675,1.020115474984113,1.054512573358387,1.039379666209665,1.023577588117426,1.053867178923646,1.022640285066947,1.056950605840609,1.025031185698543),
    @(1.02,1.019260382613656,1.051589476404707,1.020368722054587,1.054734052232546,1.039410071223072,1.02377477528621,1.053983715434875,1.022840214540975,1.057120782618311,1.025228652895725),
    @(1.02,1.019310224890516,1.051617659078566,1.020411254303518,1.054771239443004,1.039415156052785,1.023807883688424,1.054003267947927,1.022873785838922,1.057149346849043,1.025261808315667),
    @(1.02,1.018967540945175,1.05142384183394,1.020118858149709,1.054515532773243,1.03938007384309,1.02358022295466,1.053868737059498,1.022642956382733,1.056952880366419,1.025033824277548),
    @(1.02,1.017536196633623,1.050613036194153,1.018898278277446,1.053446723043209,1.039230502008633,1.022628612366196,1.053304334142525,1.021678456280069,1.056130361591622,1.024080862292897),
    @(1.02,1.015017806979407,1.049181847405989,1.016753479335922,1.051563468028275,1.038955923975571,1.020951671698537,1.052302015519072,1.01998014707541,1.054676097844083,1.022401540177142),
    @(1.02,1.013341430999438,1.048226186895506,1.015327686606929,1.050308151283393,1.038765508617955,1.019833669605001,1.051628664035733,1.018848808180997,1.053703395538302,1.021281950392398),
    @(1.02,1.012616141567497,1.047812034341163,1.014711273722165,1.049764650990227,1.038681318818633,1.019349554612805,1.051335894062308,1.018359134870692,1.053281465687463,1.020797147901383),
    @(1.02,1.012346825253086,1.047658149610549,1.014482455810613,1.049562781515395,1.038649786429092,1.019169730936508,1.051226966708844,1.018177279224385,1.053124631919415,1.02061706885484),
    @(1.02,1.012404590521573,1.047691160642257,1.014531531449805,1.049606082690091,1.038656562013371,1.019208303805968,1.051250340098693,1.018216286479913,1.053158278260299,1.020655696502198),
    @(1.02,1.012593877995138,1.04779931519136,1.014692356596622,1.049747964167706,1.038678717652743,1.019334690368043,1.051326893748042,1.01834410199633,1.053268503999708,1.020782262547689),
    @(1.02,1.012710515962773,1.047865946203748,1.014791465538436,1.049835383501721,1.038692333968045,1.019412561062913,1.051374037227345,1.018422857444642,1.053336403151522,1.020860243827877),
    @(1.02,1.013389578506626,1.04825366571964,1.015368616255075,1.050344223020027,1.03877105944258,1.019865798463051,1.051648068959276,1.018881310497995,1.053731382094956,1.021314124877111),
    @(1.02,1.013815694732875,1.048496780795881,1.015730906183931,1.050663421908296,1.038819976828493,1.020150098974774,1.051819640371348,1.019168940709178,1.053978943959276,1.02159882912817),
    @(1.02,1.014064298314519,1.048638552272506,1.015942316834394,1.050849610908542,1.038848341708608,1.02031592542053,1.051919598779284,1.0193367302052,1.054123270788886,1.021764891066499),
    @(1.02,1.014149075429106,1.048686886955305,1.0160144181781,1.050913097459328,1.038857984916946,1.020372467741212,1.051953662245449,1.019393945416212,1.054172470282622,1.021821513683762),
    @(1.02,1.013769970582162,1.04847070029397,1.015692026274668,1.050629174287689,1.038814745807787,1.020119596347193,1.05180124440134,1.019138078674629,1.053952390332086,1.021568283183359),
    @(1.02,1.012538135113958,1.047767467738732,1.014644993586232,1.049706183298946,1.038672200556236,1.019297472724874,1.051304355551686,1.018306462660031,1.053236048298221,1.020744992051198),
    @(1.02,1.011764142889329,1.047325028955632,1.013987523287937,1.049125925167812,1.038581069570026,1.018780560533841,1.050990903674331,1.017783771098216,1.052785018554863,1.020227345785563),
    @(1.02,1.012174402491723,1.047559600866767,1.014335980931573,1.049433524450729,1.038629522465378,1.019054586423381,1.051157168386374,1.018060842807159,1.053024177862541,1.020501760823297),
    @(1.02,1.013790631180733,1.048482485053661,1.015709594144484,1.050644649295839,1.038817109999423,1.020133379175341,1.051809557107801,1.01915202383878,1.053964388995245,1.021582085584704),
    @(1.02,1.015668420152855,1.04955212250763,1.017307245289623,1.052050309226438,1.039028209843364,1.021385209355299,1.052562050840683,1.020419049175509,1.055052629140668,1.022835693507105)
)

# Column order matching each inner array above.
$cols = @(2,3,4,5,6,9,10,11,12,13,14)   # B,C,D,E,F,I,J,K,L,M,N

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowVals[$j]
    }
}
